# Applies the "fixed some of the performance issues" edit:
#  - re-pairs the reference sentences (col B), the typed/input sentences (col C)
#    and the sample labels (col D) so scoring lines up correctly
#  - recalculates INTELLIGIBILITY_SCORE (col E) and Words_Correct (col G)
#  - widens column C to fit the new (longer) input-sentence text
#  - updates the last active-cell selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New reference sentences for column B (rows 2-9)
$bVals = @(
    "he is capable and willing to make decisions.",
    "Big muscles are not necessarily strong ones",
    "I think I'm getting better.",
    "You want him to do well",
    "Enjoy the fair weather while in the tropics.",
    "You're used to being on the field.",
    "We picked grapes for wine",
    "The ballet is about to begin."
)

# New input/typed sentences for column C (rows 2-9)
$cVals = @(
    "you want him to do well",
    "he is capable and willing to decide",
    "big muscles are not neccesarily strong",
    "i think i'm getting better",
    "the ballet is about to begin",
    "enjoy the fair weather while in the tropics",
    "your used to being on the feel",
    "we picked grapes for wine"
)

# New sample labels for column D (rows 2-9)
$dVals = @(
    "P1_W2_S4",
    "P1_W2_S3",
    "P1_W2_S1",
    "P1_W2_S2",
    "P1_W1_S4",
    "P1_W1_S3",
    "P1_W1_S1",
    "P1_W1_S2"
)

# Recalculated intelligibility score for column E (rows 2-9)
$eVals = @(
    0.29850746268656708,
    0.33333333333333331,
    0.2461538461538462,
    0.36734693877551022,
    0.3611111111111111,
    0.36363636363636359,
    0.32727272727272733,
    0.37037037037037029
)

# Recalculated words-correct count for column G (rows 2-9)
$gVals = @(
    8,
    7,
    6,
    8,
    8,
    8,
    7,
    8
)

for ($i = 0; $i -lt 8; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $bVals[$i]
}
for ($i = 0; $i -lt 8; $i++) {
    $row = $i + 2
    $ws.Range("C$row").Value = $cVals[$i]
}
for ($i = 0; $i -lt 8; $i++) {
    $row = $i + 2
    $ws.Range("D$row").Value = $dVals[$i]
}
for ($i = 0; $i -lt 8; $i++) {
    $row = $i + 2
    $ws.Range("E$row").Value = $eVals[$i]
}
for ($i = 0; $i -lt 8; $i++) {
    $row = $i + 2
    $ws.Range("G$row").Value = $gVals[$i]
}

# Column C is now wider to fit the new input sentences (stored width of 33)
$ws.Columns.Item(3).ColumnWidth = 32.17

# Update the last selected cell
$ws.Range("E10").Select()
